$d = $word.ActiveDocument

# Locate the "& Feiko Wielsma" portion of the subtitle run and give it
# strike-through formatting.
$r = $d.Content
$found = $r.Find.Execute("& Feiko Wielsma", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$r.Font.StrikeThrough = 1

# Word keeps a single "_GoBack" bookmark marking the last edit location;
# adding it here moves it from wherever it previously was (end of the
# "Resultaten" paragraph) to this newly-edited range.
$d.Bookmarks.Add("_GoBack", $r)
